$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-19 04:59:27"
$wsZhCn.Range("H2").Value = "2016-08-19 04:59:22"
$wsZhCn.Range("K2").Value = "2016-08-19 04:59:38"
$wsDeDe.Range("H2").Value = "2016-08-19 04:59:27"
$wsDeDe.Range("K2").Value = "2016-08-19 04:59:45"
